# The edit: the agent row that used to be "PLUG INVEST AGENTE AUTONOMO DE
# INVESTIMENTOS S/S LTDA" (row 226) was renamed to "SOMMA ASSESSORIA DE
# INVESTIMENTOS S/S LTDA" (all other columns for that record stay the same).
# Because the sheet keeps rows ordered by company name, this rename moves the
# record down past all the other rows whose name starts with "PR..", "PO..",
# "Pr..", "R..", "S.." up to "SLR ..." (rows 227-268), which therefore all
# shift up by one row. The net effect on rows 226-268 is:
#   after[226..267] = before[227..268]
#   after[268]       = before[226], but with the company name (column A)
#                       changed to "SOMMA ASSESSORIA DE INVESTIMENTOS S/S LTDA"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 226
$lastRow  = 268
$newName  = "SOMMA ASSESSORIA DE INVESTIMENTOS S/S LTDA"

# --- Capture all the data we need BEFORE writing anything -------------------

# Original content of the row that is being renamed (kept for columns B..L).
$origFirstRow = $ws.Range("A" + $firstRow + ":L" + $firstRow).Value()

# Original content of all the rows that will shift up by one.
$block = $ws.Range("A" + ($firstRow + 1) + ":L" + $lastRow).Value()

# --- Write the shifted block into rows 226..267 ------------------------------

$destBlock = $ws.Range("A" + $firstRow + ":L" + ($lastRow - 1))
# Force text formatting while assigning so that values which look like dates
# (e.g. "07/02/2022") are not silently converted into real Excel dates.
$destBlock.NumberFormat = "@"
$destBlock.Value = $block
$destBlock.NumberFormat = "General"

# --- Write the renamed row into row 268 --------------------------------------

$lastRowRange = $ws.Range("A" + $lastRow + ":L" + $lastRow)
$lastRowRange.NumberFormat = "@"

$lastRowRange.Value = $origFirstRow
$ws.Range("A" + $lastRow).Value = $newName

$lastRowRange.NumberFormat = "General"
